# Update gh-pages to output generated at 456a3b4
# Refreshes the "想去人数" (F column) figures across all four sheets.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 4581
$ws1.Cells.Item(5, 6).Value = 3690
$ws1.Cells.Item(6, 6).Value = 1073
$ws1.Cells.Item(9, 6).Value = 372
$ws1.Cells.Item(10, 6).Value = 370
$ws1.Cells.Item(11, 6).Value = 2562
$ws1.Cells.Item(14, 6).Value = 1968
$ws1.Cells.Item(15, 6).Value = 277
$ws1.Cells.Item(17, 6).Value = 559
$ws1.Cells.Item(20, 6).Value = 10608
$ws1.Cells.Item(21, 6).Value = 6137
$ws1.Cells.Item(30, 6).Value = 187
$ws1.Cells.Item(31, 6).Value = 865
$ws1.Cells.Item(32, 6).Value = 3570
$ws1.Cells.Item(34, 6).Value = 971
$ws1.Cells.Item(36, 6).Value = 134
$ws1.Cells.Item(37, 6).Value = 279
$ws1.Cells.Item(39, 6).Value = 255
$ws1.Cells.Item(44, 6).Value = 203
$ws1.Cells.Item(45, 6).Value = 118

# --- 演出 (sheet 2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(15, 6).Value = 3608

# --- 本地生活 (sheet 3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 8852
$ws3.Cells.Item(4, 6).Value = 1672

# --- 全部类型 (sheet 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 8852
$ws4.Cells.Item(4, 6).Value = 1672
$ws4.Cells.Item(5, 6).Value = 4581
$ws4.Cells.Item(7, 6).Value = 3690
$ws4.Cells.Item(8, 6).Value = 1073
$ws4.Cells.Item(11, 6).Value = 370
$ws4.Cells.Item(12, 6).Value = 2562
$ws4.Cells.Item(19, 6).Value = 277
$ws4.Cells.Item(21, 6).Value = 559
$ws4.Cells.Item(23, 6).Value = 10608
$ws4.Cells.Item(24, 6).Value = 3608
$ws4.Cells.Item(33, 6).Value = 865
$ws4.Cells.Item(34, 6).Value = 3570
$ws4.Cells.Item(36, 6).Value = 971
$ws4.Cells.Item(37, 6).Value = 134
$ws4.Cells.Item(38, 6).Value = 279
$ws4.Cells.Item(41, 6).Value = 255
$ws4.Cells.Item(46, 6).Value = 118
